$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "ECTS" column header in K1, matching the formatting of the
# other header cells (B1:I1).
$ws.Range("K1").Value = "ECTS"
$ws.Range("I1").Copy()
$ws.Range("K1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in ECTS credit values for each module row (2-24). Every module is
# worth 6 ECTS, except the Master Thesis (row 22), which is worth 18.
$ws.Range("K2:K21").Value = 6
$ws.Range("K22").Value = 18
$ws.Range("K23:K24").Value = 6

# Reflect the author's final view/selection state: scrolled right with
# K24 as the active cell.
$ws.Range("K24").Select() | Out-Null
